$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows at 332-334 (shifts existing rows 332.. down by 3)
$ws.Range("A332:A334").EntireRow.Insert()

# Row 332 - new data row (Choclero, Primera)
$ws.Range("A332").Value = 3
$ws.Range("B332").Value = "Femacal de La Calera"
$ws.Range("C332").Value = "Coquimbo"
$ws.Range("D332").Value = 44543
$ws.Range("E332").Value = 5
$ws.Range("F332").Value = 100112024
$ws.Range("G332").Value = "Choclo"
$ws.Range("H332").Value = "Choclero"
$ws.Range("I332").Value = "Primera"
$ws.Range("J332").Value = 12000
$ws.Range("K332").Value = 550
$ws.Range("L332").Value = 600
$ws.Range("M332").Value = 574
$ws.Range("N332").Value = "$/unidad"
$ws.Range("O332").Value = "Provincia de Limarí"
$ws.Range("P332").Value = 574
$ws.Range("Q332").Value = 1
$ws.Range("R332").Value = "Hortaliza"

# Row 333 - new data row (Choclero, Segunda)
$ws.Range("A333").Value = 3
$ws.Range("B333").Value = "Femacal de La Calera"
$ws.Range("C333").Value = "Coquimbo"
$ws.Range("D333").Value = 44543
$ws.Range("E333").Value = 5
$ws.Range("F333").Value = 100112024
$ws.Range("G333").Value = "Choclo"
$ws.Range("H333").Value = "Choclero"
$ws.Range("I333").Value = "Segunda"
$ws.Range("J333").Value = 5500
$ws.Range("K333").Value = 400
$ws.Range("L333").Value = 400
$ws.Range("M333").Value = 400
$ws.Range("N333").Value = "$/unidad"
$ws.Range("O333").Value = "Provincia de Limarí"
$ws.Range("P333").Value = 400
$ws.Range("Q333").Value = 1
$ws.Range("R333").Value = "Hortaliza"

# Row 334 - new data row (Dulce o Americano, Primera)
$ws.Range("A334").Value = 3
$ws.Range("B334").Value = "Femacal de La Calera"
$ws.Range("C334").Value = "Coquimbo"
$ws.Range("D334").Value = 44543
$ws.Range("E334").Value = 5
$ws.Range("F334").Value = 100112024
$ws.Range("G334").Value = "Choclo"
$ws.Range("H334").Value = "Dulce o Americano"
$ws.Range("I334").Value = "Primera"
$ws.Range("J334").Value = 2800
$ws.Range("K334").Value = 300
$ws.Range("L334").Value = 300
$ws.Range("M334").Value = 300
$ws.Range("N334").Value = "$/unidad"
$ws.Range("O334").Value = "Provincia de Limarí"
$ws.Range("P334").Value = 300
$ws.Range("Q334").Value = 1
$ws.Range("R334").Value = "Hortaliza"
